$wb = $excel.ActiveWorkbook


# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 654398.5
$ws.Range("J17").Value = 654398.5
$ws.Range("L17").Value = 1963195.5
$ws.Range("N17").Value = -1963531.5
$ws.Range("H116").Value = 2440.2856
$ws.Range("I116").Value = 3017.5
$ws.Range("J116").Value = 2007.375
$ws.Range("K116").Value = 3017.5
$ws.Range("L116").Value = 2007.375
$ws.Range("M116").Value = 424.5
$ws.Range("N116").Value = -8891.375
$ws.Range("H137").Value = 3699.9592
$ws.Range("I137").Value = 1335.4286
$ws.Range("J137").Value = 6852.6665
$ws.Range("K137").Value = 4006.2858
$ws.Range("L137").Value = 20557.9995
$ws.Range("M137").Value = -1456.2858
$ws.Range("N137").Value = -25657.9995

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 111967.336
$ws.Range("I2").Value = 111967.336
$ws.Range("K2").Value = 111967.336
$ws.Range("M2").Value = -111854.336
$ws.Range("H32").Value = 7091.041
$ws.Range("I32").Value = 5991.659
$ws.Range("K32").Value = 5991.659
$ws.Range("M32").Value = -5704.659
$ws.Range("H45").Value = 9135.1
$ws.Range("I45").Value = 12183.857
$ws.Range("K45").Value = 12183.857
$ws.Range("M45").Value = -11806.857
$ws.Range("H74").Value = 4493.931
$ws.Range("I74").Value = 1021.56525
$ws.Range("J74").Value = 17804.666
$ws.Range("K74").Value = 1021.56525
$ws.Range("L74").Value = 17804.666
$ws.Range("M74").Value = -147.56525
$ws.Range("N74").Value = -19552.666
$ws.Range("H77").Value = 4493.931
$ws.Range("I77").Value = 1021.56525
$ws.Range("J77").Value = 17804.666
$ws.Range("K77").Value = 5107.82625
$ws.Range("L77").Value = 89023.33
$ws.Range("M77").Value = -739.8262500000001
$ws.Range("N77").Value = -97759.33
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H116").Value = 111967.336
$ws.Range("I116").Value = 111967.336
$ws.Range("K116").Value = 111967.336
$ws.Range("M116").Value = -109673.336
$ws.Range("H122").Value = 2311.7693
$ws.Range("I122").Value = 2326.5557
$ws.Range("J122").Value = 2278.5
$ws.Range("K122").Value = 6979.6671
$ws.Range("L122").Value = 6835.5
$ws.Range("M122").Value = -4529.6671
$ws.Range("N122").Value = -11735.5
$ws.Range("H132").Value = 12704.95
$ws.Range("I132").Value = 13738.777
$ws.Range("J132").Value = 11859.091
$ws.Range("K132").Value = 41216.331
$ws.Range("L132").Value = 35577.273
$ws.Range("M132").Value = -38686.331
$ws.Range("N132").Value = -40637.273

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 111967.336
$ws.Range("I3").Value = 111967.336
$ws.Range("K3").Value = 111967.336
$ws.Range("M3").Value = -111853.336
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H105").Value = 1703.5333
$ws.Range("I105").Value = 1285
$ws.Range("J105").Value = 1855.7273
$ws.Range("K105").Value = 1285
$ws.Range("L105").Value = 1855.7273
$ws.Range("M105").Value = 462
$ws.Range("N105").Value = -5349.7273

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 16198.25
$ws.Range("I37").Value = 4622
$ws.Range("K37").Value = 4622
$ws.Range("M37").Value = -4515
$ws.Range("H111").Value = 32000
$ws.Range("J111").Value = 32000
$ws.Range("L111").Value = 32000
$ws.Range("N111").Value = -40180
$ws.Range("H132").Value = 60613720
$ws.Range("I132").Value = 88898264
$ws.Range("J132").Value = 3989.4285
$ws.Range("K132").Value = 266694792
$ws.Range("L132").Value = 11968.2855
$ws.Range("M132").Value = -266692262
$ws.Range("N132").Value = -17028.2855
$ws.Range("H134").Value = 1981.9565
$ws.Range("I134").Value = 2069.25
$ws.Range("J134").Value = 1400
$ws.Range("K134").Value = 6207.75
$ws.Range("L134").Value = 4200
$ws.Range("M134").Value = -3672.75
$ws.Range("N134").Value = -9270

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 789.1667
$ws.Range("I5").Value = 410
$ws.Range("J5").Value = 888.9474
$ws.Range("K5").Value = 1230
$ws.Range("L5").Value = 2666.8422
$ws.Range("M5").Value = -1118
$ws.Range("N5").Value = -2890.8422
$ws.Range("H12").Value = 31.428572
$ws.Range("J12").Value = 43.15
$ws.Range("L12").Value = 129.45
$ws.Range("N12").Value = -475.45
$ws.Range("H122").Value = 1068.75
$ws.Range("I122").Value = 402.44446
$ws.Range("J122").Value = 1329.4783
$ws.Range("K122").Value = 3622.00014
$ws.Range("L122").Value = 11965.3047
$ws.Range("M122").Value = -1172.00014
$ws.Range("N122").Value = -16865.3047
$ws.Range("H135").Value = 789.1667
$ws.Range("I135").Value = 410
$ws.Range("J135").Value = 888.9474
$ws.Range("K135").Value = 3690
$ws.Range("L135").Value = 8000.5266
$ws.Range("M135").Value = -1155
$ws.Range("N135").Value = -13070.5266
$ws.Range("H141").Value = 146508.42
$ws.Range("I141").Value = 146508.42
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 439525.26
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -434345.26
$ws.Range("N141").ClearContents()

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1977.75
$ws.Range("J7").Value = 1400
$ws.Range("L7").Value = 1400
$ws.Range("N7").Value = -1624
$ws.Range("H8").Value = 1977.75
$ws.Range("J8").Value = 1400
$ws.Range("L8").Value = 1400
$ws.Range("N8").Value = -1678

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 1050
$ws.Range("I5").Value = 1200
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -1087
$ws.Range("N5").Value = -1226
$ws.Range("H132").Value = 8664.611000000001
$ws.Range("I132").Value = 10004.643
$ws.Range("J132").Value = 3974.5
$ws.Range("K132").Value = 30013.929
$ws.Range("L132").Value = 11923.5
$ws.Range("M132").Value = -27483.929
$ws.Range("N132").Value = -16983.5

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1054855.8
$ws.Range("I14").Value = 2501445
$ws.Range("J14").Value = 2790.818
$ws.Range("K14").Value = 2501445
$ws.Range("L14").Value = 2790.818
$ws.Range("M14").Value = -2501277
$ws.Range("N14").Value = -3126.818
